$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (42 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 2747.25  # H4: 1612.5714 -> 2747.25
$ws.Cells.Item(4, 9).Value = 3533  # I4: 1816.3334 -> 3533
$ws.Cells.Item(4, 11).Value = 3533  # K4: 1816.3334 -> 3533
$ws.Cells.Item(4, 13).Value = -3419  # M4: -1702.3334 -> -3419
$ws.Cells.Item(11, 8).Value = 49.235294  # H11: 46.555557 -> 49.235294
$ws.Cells.Item(11, 9).Value = 49.235294  # I11: 46.555557 -> 49.235294
$ws.Cells.Item(11, 11).Value = 49.235294  # K11: 46.555557 -> 49.235294
$ws.Cells.Item(11, 13).Value = 90.76470599999999  # M11: 93.44444300000001 -> 90.76470599999999
$ws.Cells.Item(17, 8).Value = 5798.077  # H17: 5256.6 -> 5798.077
$ws.Cells.Item(17, 9).Value = 1500  # I17: 1499.3334 -> 1500
$ws.Cells.Item(17, 10).Value = 5970  # J17: 5674.074 -> 5970
$ws.Cells.Item(17, 11).Value = 4500  # K17: 4498.0002 -> 4500
$ws.Cells.Item(17, 12).Value = 17910  # L17: 17022.222 -> 17910
$ws.Cells.Item(17, 13).Value = -4332  # M17: -4330.0002 -> -4332
$ws.Cells.Item(17, 14).Value = -18246  # N17: -17358.222 -> -18246
$ws.Cells.Item(58, 8).Value = 2175  # H58: 2147.625 -> 2175
$ws.Cells.Item(58, 9).Value = 1135  # I58: 946.8333 -> 1135
$ws.Cells.Item(58, 10).Value = 3475  # J58: 5750 -> 3475
$ws.Cells.Item(58, 11).Value = 3405  # K58: 2840.4999 -> 3405
$ws.Cells.Item(58, 12).Value = 10425  # L58: 17250 -> 10425
$ws.Cells.Item(58, 13).Value = -3255  # M58: -2690.4999 -> -3255
$ws.Cells.Item(58, 14).Value = -10725  # N58: -17550 -> -10725
$ws.Cells.Item(64, 8).Value = 5250  # H64: 0 -> 5250
$ws.Cells.Item(64, 9).Value = 5000  # I64: 0 -> 5000
$ws.Cells.Item(64, 10).Value = 5500  # J64: 0 -> 5500
$ws.Cells.Item(64, 11).Value = 5000  # K64: 0 -> 5000
$ws.Cells.Item(64, 12).Value = 5500  # L64: 0 -> 5500
$ws.Cells.Item(64, 13).Value = -4752  # M64: None -> -4752
$ws.Cells.Item(64, 14).Value = -5996  # N64: None -> -5996
$ws.Cells.Item(67, 8).Value = 5250  # H67: 0 -> 5250
$ws.Cells.Item(67, 9).Value = 5000  # I67: 0 -> 5000
$ws.Cells.Item(67, 10).Value = 5500  # J67: 0 -> 5500
$ws.Cells.Item(67, 11).Value = 5000  # K67: 0 -> 5000
$ws.Cells.Item(67, 12).Value = 5500  # L67: 0 -> 5500
$ws.Cells.Item(67, 13).Value = -4142  # M67: None -> -4142
$ws.Cells.Item(67, 14).Value = -7216  # N67: None -> -7216
$ws.Cells.Item(76, 8).Value = 3600  # H76: 3300 -> 3600
$ws.Cells.Item(79, 8).Value = 3600  # H79: 3300 -> 3600
$ws.Cells.Item(137, 8).Value = 6382.8335  # H137: 7259.6 -> 6382.8335
$ws.Cells.Item(137, 9).Value = 5959.4  # I137: 6949.5 -> 5959.4
$ws.Cells.Item(137, 11).Value = 17878.2  # K137: 20848.5 -> 17878.2
$ws.Cells.Item(137, 13).Value = -15328.2  # M137: -18298.5 -> -15328.2

# --- Sheet: ARM (16 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(23, 8).Value = 21666.334  # H23: 14999.6 -> 21666.334
$ws.Cells.Item(23, 9).Value = 5000  # I23: 4999.6665 -> 5000
$ws.Cells.Item(23, 11).Value = 5000  # K23: 4999.6665 -> 5000
$ws.Cells.Item(23, 13).Value = -4741  # M23: -4740.6665 -> -4741
$ws.Cells.Item(37, 8).Value = 2518.3333  # H37: 3500 -> 2518.3333
$ws.Cells.Item(37, 10).Value = 555  # J37: 0 -> 555
$ws.Cells.Item(37, 12).Value = 555  # L37: 0 -> 555
$ws.Cells.Item(37, 14).Value = -1101  # N37: None -> -1101
$ws.Cells.Item(45, 8).Value = 0  # H45: 400 -> 0
$ws.Cells.Item(45, 9).Value = 0  # I45: 400 -> 0
$ws.Cells.Item(45, 11).Value = 0  # K45: 400 -> 0
$ws.Cells.Item(45, 13).Value = $null  # M45: -23 -> (cleared)
$ws.Cells.Item(127, 8).Value = 92856  # H127: 52999 -> 92856
$ws.Cells.Item(127, 10).Value = 92856  # J127: 52999 -> 92856
$ws.Cells.Item(127, 12).Value = 92856  # L127: 52999 -> 92856
$ws.Cells.Item(127, 14).Value = -102776  # N127: -62919 -> -102776

# --- Sheet: BSM (7 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(14, 8).Value = 0  # H14: 549 -> 0
$ws.Cells.Item(14, 9).Value = 0  # I14: 999 -> 0
$ws.Cells.Item(14, 10).Value = 0  # J14: 99 -> 0
$ws.Cells.Item(14, 11).Value = 0  # K14: 999 -> 0
$ws.Cells.Item(14, 12).Value = 0  # L14: 99 -> 0
$ws.Cells.Item(14, 13).Value = $null  # M14: -827 -> (cleared)
$ws.Cells.Item(14, 14).Value = $null  # N14: -443 -> (cleared)

# --- Sheet: CRP (36 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 16464.545  # H31: 15259.083 -> 16464.545
$ws.Cells.Item(31, 10).Value = 9199.799999999999  # J31: 7999.6665 -> 9199.799999999999
$ws.Cells.Item(31, 12).Value = 9199.799999999999  # L31: 7999.6665 -> 9199.799999999999
$ws.Cells.Item(31, 14).Value = -9789.799999999999  # N31: -8589.666499999999 -> -9789.799999999999
$ws.Cells.Item(34, 8).Value = 16464.545  # H34: 15259.083 -> 16464.545
$ws.Cells.Item(34, 10).Value = 9199.799999999999  # J34: 7999.6665 -> 9199.799999999999
$ws.Cells.Item(34, 12).Value = 9199.799999999999  # L34: 7999.6665 -> 9199.799999999999
$ws.Cells.Item(34, 14).Value = -9603.799999999999  # N34: -8403.666499999999 -> -9603.799999999999
$ws.Cells.Item(62, 8).Value = 3234.1667  # H62: 3661 -> 3234.1667
$ws.Cells.Item(62, 9).Value = 1735  # I62: 2052.5 -> 1735
$ws.Cells.Item(62, 11).Value = 1735  # K62: 2052.5 -> 1735
$ws.Cells.Item(62, 13).Value = -1111  # M62: -1428.5 -> -1111
$ws.Cells.Item(65, 8).Value = 3234.1667  # H65: 3661 -> 3234.1667
$ws.Cells.Item(65, 9).Value = 1735  # I65: 2052.5 -> 1735
$ws.Cells.Item(65, 11).Value = 8675  # K65: 10262.5 -> 8675
$ws.Cells.Item(65, 13).Value = -5555  # M65: -7142.5 -> -5555
$ws.Cells.Item(74, 8).Value = 52249  # H74: 52249.5 -> 52249
$ws.Cells.Item(74, 9).Value = 9498  # I74: 9499 -> 9498
$ws.Cells.Item(74, 11).Value = 9498  # K74: 9499 -> 9498
$ws.Cells.Item(74, 13).Value = -8624  # M74: -8625 -> -8624
$ws.Cells.Item(77, 8).Value = 52249  # H77: 52249.5 -> 52249
$ws.Cells.Item(77, 9).Value = 9498  # I77: 9499 -> 9498
$ws.Cells.Item(77, 11).Value = 28494  # K77: 28497 -> 28494
$ws.Cells.Item(77, 13).Value = -24126  # M77: -24129 -> -24126
$ws.Cells.Item(99, 8).Value = 42446.89  # H99: 47252.75 -> 42446.89
$ws.Cells.Item(99, 9).Value = 52337.168  # I99: 62004.6 -> 52337.168
$ws.Cells.Item(99, 11).Value = 52337.168  # K99: 62004.6 -> 52337.168
$ws.Cells.Item(99, 13).Value = -50839.168  # M99: -60506.6 -> -50839.168
$ws.Cells.Item(126, 8).Value = 42446.89  # H126: 47252.75 -> 42446.89
$ws.Cells.Item(126, 9).Value = 52337.168  # I126: 62004.6 -> 52337.168
$ws.Cells.Item(126, 11).Value = 157011.504  # K126: 186013.8 -> 157011.504
$ws.Cells.Item(126, 13).Value = -154541.504  # M126: -183543.8 -> -154541.504
$ws.Cells.Item(134, 8).Value = 3645.45  # H134: 3605.1904 -> 3645.45
$ws.Cells.Item(134, 9).Value = 3618.1765  # I134: 3572.7222 -> 3618.1765
$ws.Cells.Item(134, 11).Value = 10854.5295  # K134: 10718.1666 -> 10854.5295
$ws.Cells.Item(134, 13).Value = -8319.529500000001  # M134: -8183.1666 -> -8319.529500000001

# --- Sheet: CUL (45 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 109.5  # H26: 138.42857 -> 109.5
$ws.Cells.Item(26, 9).Value = 87.72727  # I26: 54 -> 87.72727
$ws.Cells.Item(26, 10).Value = 349  # J26: 349.5 -> 349
$ws.Cells.Item(26, 11).Value = 263.18181  # K26: 162 -> 263.18181
$ws.Cells.Item(26, 12).Value = 1047  # L26: 1048.5 -> 1047
$ws.Cells.Item(26, 13).Value = 24.81818999999996  # M26: 126 -> 24.81818999999996
$ws.Cells.Item(26, 14).Value = -1623  # N26: -1624.5 -> -1623
$ws.Cells.Item(34, 8).Value = 3493.3333  # H34: 4095.3 -> 3493.3333
$ws.Cells.Item(34, 10).Value = 3810.7273  # J34: 4550.1113 -> 3810.7273
$ws.Cells.Item(34, 12).Value = 11432.1819  # L34: 13650.3339 -> 11432.1819
$ws.Cells.Item(34, 14).Value = -11600.1819  # N34: -13818.3339 -> -11600.1819
$ws.Cells.Item(55, 8).Value = 1806.1428  # H55: 1893.2 -> 1806.1428
$ws.Cells.Item(55, 10).Value = 2044.1818  # J55: 2133.1667 -> 2044.1818
$ws.Cells.Item(55, 12).Value = 6132.5454  # L55: 6399.500100000001 -> 6132.5454
$ws.Cells.Item(55, 14).Value = -6486.5454  # N55: -6753.500100000001 -> -6486.5454
$ws.Cells.Item(86, 8).Value = 290  # H86: 291.33334 -> 290
$ws.Cells.Item(86, 9).Value = 293.33334  # I86: 291.33334 -> 293.33334
$ws.Cells.Item(86, 10).Value = 280  # J86: 0 -> 280
$ws.Cells.Item(86, 11).Value = 880.0000200000001  # K86: 874.0000200000001 -> 880.0000200000001
$ws.Cells.Item(86, 12).Value = 840  # L86: 0 -> 840
$ws.Cells.Item(86, 13).Value = 305.9999799999999  # M86: 311.9999799999999 -> 305.9999799999999
$ws.Cells.Item(86, 14).Value = -3212  # N86: None -> -3212
$ws.Cells.Item(89, 8).Value = 290  # H89: 291.33334 -> 290
$ws.Cells.Item(89, 9).Value = 293.33334  # I89: 291.33334 -> 293.33334
$ws.Cells.Item(89, 10).Value = 280  # J89: 0 -> 280
$ws.Cells.Item(89, 11).Value = 2640.00006  # K89: 2622.00006 -> 2640.00006
$ws.Cells.Item(89, 12).Value = 2520  # L89: 0 -> 2520
$ws.Cells.Item(89, 13).Value = 3287.99994  # M89: 3305.99994 -> 3287.99994
$ws.Cells.Item(89, 14).Value = -14376  # N89: None -> -14376
$ws.Cells.Item(100, 8).Value = 5000  # H100: 0 -> 5000
$ws.Cells.Item(100, 10).Value = 5000  # J100: 0 -> 5000
$ws.Cells.Item(100, 12).Value = 15000  # L100: 0 -> 15000
$ws.Cells.Item(100, 14).Value = -16622  # N100: None -> -16622
$ws.Cells.Item(113, 8).Value = 2642.3157  # H113: 2693.5 -> 2642.3157
$ws.Cells.Item(113, 10).Value = 2642  # J113: 2701.5833 -> 2642
$ws.Cells.Item(113, 12).Value = 7926  # L113: 8104.749899999999 -> 7926
$ws.Cells.Item(113, 14).Value = -12266  # N113: -12444.7499 -> -12266
$ws.Cells.Item(131, 8).Value = 3066.4443  # H131: 2088.25 -> 3066.4443
$ws.Cells.Item(131, 9).Value = 0  # I131: 830.5714 -> 0
$ws.Cells.Item(131, 11).Value = 0  # K131: 2491.7142 -> 0
$ws.Cells.Item(131, 13).Value = $null  # M131: 2548.2858 -> (cleared)
$ws.Cells.Item(134, 8).Value = 1344.8  # H134: 1544.8 -> 1344.8
$ws.Cells.Item(134, 9).Value = 1344.8  # I134: 1544.8 -> 1344.8
$ws.Cells.Item(134, 11).Value = 4034.4  # K134: 4634.4 -> 4034.4
$ws.Cells.Item(134, 13).Value = 1035.6  # M134: 435.6000000000004 -> 1035.6

# --- Sheet: GSM (11 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 264.46155  # H2: 263.76923 -> 264.46155
$ws.Cells.Item(2, 9).Value = 355.66666  # I2: 324.1 -> 355.66666
$ws.Cells.Item(2, 10).Value = 59.25  # J2: 62.666668 -> 59.25
$ws.Cells.Item(2, 11).Value = 355.66666  # K2: 324.1 -> 355.66666
$ws.Cells.Item(2, 12).Value = 59.25  # L2: 62.666668 -> 59.25
$ws.Cells.Item(2, 13).Value = -242.66666  # M2: -211.1 -> -242.66666
$ws.Cells.Item(2, 14).Value = -285.25  # N2: -288.666668 -> -285.25
$ws.Cells.Item(107, 8).Value = 532.25  # H107: 532.3333 -> 532.25
$ws.Cells.Item(107, 9).Value = 491.54544  # I107: 491.63635 -> 491.54544
$ws.Cells.Item(107, 11).Value = 491.54544  # K107: 491.63635 -> 491.54544
$ws.Cells.Item(107, 13).Value = 1428.45456  # M107: 1428.36365 -> 1428.45456

# --- Sheet: LTW (19 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1745.9131  # H22: 1950.6 -> 1745.9131
$ws.Cells.Item(22, 9).Value = 1619.3684  # I22: 1851.5 -> 1619.3684
$ws.Cells.Item(22, 11).Value = 1619.3684  # K22: 1851.5 -> 1619.3684
$ws.Cells.Item(22, 13).Value = -1324.3684  # M22: -1556.5 -> -1324.3684
$ws.Cells.Item(27, 8).Value = 1745.9131  # H27: 1950.6 -> 1745.9131
$ws.Cells.Item(27, 9).Value = 1619.3684  # I27: 1851.5 -> 1619.3684
$ws.Cells.Item(27, 11).Value = 1619.3684  # K27: 1851.5 -> 1619.3684
$ws.Cells.Item(27, 13).Value = -1512.3684  # M27: -1744.5 -> -1512.3684
$ws.Cells.Item(55, 8).Value = 3885  # H55: 2454.6667 -> 3885
$ws.Cells.Item(55, 9).Value = 5788  # I55: 2130 -> 5788
$ws.Cells.Item(55, 10).Value = 3409.25  # J55: 2779.3333 -> 3409.25
$ws.Cells.Item(55, 11).Value = 5788  # K55: 2130 -> 5788
$ws.Cells.Item(55, 12).Value = 3409.25  # L55: 2779.3333 -> 3409.25
$ws.Cells.Item(55, 13).Value = -5615  # M55: -1957 -> -5615
$ws.Cells.Item(55, 14).Value = -3755.25  # N55: -3125.3333 -> -3755.25
$ws.Cells.Item(132, 8).Value = 5473  # H132: 5420.974 -> 5473
$ws.Cells.Item(132, 9).Value = 4504.9165  # I132: 4462.48 -> 4504.9165
$ws.Cells.Item(132, 11).Value = 13514.7495  # K132: 13387.44 -> 13514.7495
$ws.Cells.Item(132, 13).Value = -10984.7495  # M132: -10857.44 -> -10984.7495

# --- Sheet: WVR (8 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 5000  # H122: 3833.3333 -> 5000
$ws.Cells.Item(122, 9).Value = 0  # I122: 1500 -> 0
$ws.Cells.Item(122, 11).Value = 0  # K122: 4500 -> 0
$ws.Cells.Item(122, 13).Value = $null  # M122: -2050 -> (cleared)
$ws.Cells.Item(132, 8).Value = 1877.6  # H132: 2166.3333 -> 1877.6
$ws.Cells.Item(132, 9).Value = 1972  # I132: 2499.5 -> 1972
$ws.Cells.Item(132, 11).Value = 5916  # K132: 7498.5 -> 5916
$ws.Cells.Item(132, 13).Value = -3386  # M132: -4968.5 -> -3386
